# CLSD.xlsx update: add CLS-AX sheet, cross-link with Main, update price/multiples.
$wb = $excel.ActiveWorkbook
$main = $wb.Worksheets.Item("Main")

# --- Add the new "CLS-AX" worksheet right after "Main" ---
$clsax = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $main)
$clsax.Name = "CLS-AX"

# Link back to Main from the new sheet.
$clsax.Hyperlinks.Add($clsax.Range("A1"), "", "Main!A1", "", "Main")

# Brand / generic / trial summary table on CLS-AX.
$clsax.Range("B2").Value = "Brand"
$clsax.Range("C2").Value = "CLS-AX"

# Update the hyperlink on Main!B3 to point at the new sheet.
$main.Hyperlinks.Add($main.Range("B3"), "", "'CLS-AX'!A1", "", "CLS-AX (axitinib)")

$clsax.Range("B3").Value = "Generic"
$clsax.Range("C3").Value = "axitinib"

$clsax.Range("B4").Value = "Clinical Trials"

$clsax.Range("C5").Value = 'Phase IIb "ODYSSEY" n=60 wet AMD 36 weeks - NCT05891548'
$clsax.Range("C5").Font.Bold = $true
$clsax.Range("C5").Font.Underline = $true

$clsax.Range("C10").Value = 'Phase I/II "OASIS"'
$clsax.Range("C10").Font.Bold = $true
$clsax.Range("C10").Font.Underline = $true

$clsax.Range("C6").Value = "CLS-AX vs aflibercept (2:1)"

$clsax.Range("C7").Value = "PE: BCVA from baseline"
$clsax.Range("C7").Characters(15, 8).Font.Bold = $true

$clsax.Range("C11").Value = "Results : Four dose-escalating cohorts (0.03 mg, n=6; 0.1 mg, n=5; 0.5 mg, n=8; 1.0 mg, n=8) were enrolled with a mean age 81 years, mean duration of nAMD diagnosis 54 months, and mean 29.9 prior anti-VEGF injections. In all cohorts, there were no serious adverse events, no treatment emergent adverse events related to study treatment, no dose limiting toxicities, and no adverse events related to inflammation, vasculitis or vascular occlusion. In cohorts 3 and 4, in OASIS to the 3-month timepoint (n=16), there was a 73% reduction in treatment burden, and in the ongoing Extension Study (interim data as of 10/27/22, n=12), there was a 90% reduction in treatment burden from the average monthly injections before CLS-AX. In cohorts 3 and 4 from the ongoing Extension Study, the injection free rate for supplemental aflibercept treatment was 88% up to month 5 (7/8 patients not receiving additional therapy) and 75% to month 6 (3/4 of patients not receiving additional therapy). Through 6 months, stable mean BCVA and anatomic signs of biological effect, including stable mean CST, were observed."

# --- Column widths / layout on CLS-AX (best-fit column widths) ---
$clsax.Columns.Item(1).ColumnWidth = 4.17
$clsax.Columns.Item(2).ColumnWidth = 11.17

# --- Column widths / layout on Main ---
$main.Columns.Item(1).ColumnWidth = 3.34
$main.Columns.Item(2).ColumnWidth = 15.65
$main.Range("D1:H1").ColumnWidth = 6.5

# --- Price update on Main (flows through the EV build formulas) ---
$main.Range("K2").Value = 1.27

# --- View state: CLS-AX becomes the active / visible tab, zoomed differently than Main ---
$main.Range("E11").Select()
$main.Application.ActiveWindow.Zoom = 160

$clsax.Activate()
$clsax.Range("E9").Select()
$clsax.Application.ActiveWindow.Zoom = 250
